$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value  = 12.5311
$ws.Range("C6").Value  = -11.73549999999999
$ws.Range("C7").Value  = -11.90769999999999
$ws.Range("B8").Value  = 5.199899999999997
$ws.Range("E8").Value  = 13.08629999999999
$ws.Range("A12").Value = -22.74950000000001
$ws.Range("B12").Value = 5.739000000000001
$ws.Range("B14").Value = 8.453100000000003
$ws.Range("C19").Value = -13.08569999999999
$ws.Range("D19").Value = -8.194499999999994
$ws.Range("E19").Value = 13.476
$ws.Range("C21").Value = -12.8351
$ws.Range("B22").Value = 5.014200000000006
$ws.Range("C24").Value = -11.4084
